$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text-preserving number format on cells whose new values would
# otherwise be auto-converted to numbers (losing significant trailing zeros)
# by the COM Range.Value setter, matching the original text-typed cells.
foreach ($addr in @('D10','D25','D26','D36','D47')) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '72.231.86'
$ws.Range('E2').Value = '  -0.17%  '
$ws.Range('D3').Value = '2.652.60'
$ws.Range('E3').Value = '  +1.24%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '597.14'
$ws.Range('E5').Value = '  -1.08%  '
$ws.Range('D6').Value = '174.84'
$ws.Range('E6').Value = '  -2.09%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').Value = '0.524'
$ws.Range('E8').Value = '  -0.34%  '
$ws.Range('D9').Value = '2.650.05'
$ws.Range('E9').Value = '  +1.19%  '
$ws.Range('D10').Value = '0.170'
$ws.Range('E10').Value = '  -2.15%  '
$ws.Range('E11').Value = '  +2.10%  '
$ws.Range('E12').Value = '  +0.91%  '
$ws.Range('E13').Value = '  -0.82%  '
$ws.Range('D14').Value = '3.135.87'
$ws.Range('E14').Value = '  +1.07%  '
$ws.Range('B15').Value = 'ShibaInu'
$ws.Range('C15').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D15').Value = '0.0000185'
$ws.Range('E15').Value = '  -1.48%  '
$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').Value = '72.117.32'
$ws.Range('E16').Value = '  -0.24%  '
$ws.Range('D17').Value = '26.23'
$ws.Range('E17').Value = '  -1.55%  '
$ws.Range('D18').Value = '2.644.58'
$ws.Range('E18').Value = '  +0.81%  '
$ws.Range('D19').Value = '12.22'
$ws.Range('E19').Value = '  +5.37%  '
$ws.Range('D20').Value = '8.14'
$ws.Range('E20').Value = '  +3.54%  '
$ws.Range('D21').Value = '370.36'
$ws.Range('E21').Value = '  -2.65%  '
$ws.Range('D22').Value = '4.18'
$ws.Range('E22').Value = '  -0.04%  '
$ws.Range('E23').Value = '  +0.04%  '
$ws.Range('D24').Value = '72.08'
$ws.Range('E24').Value = '  -1.77%  '
$ws.Range('D25').Value = '1.00'
$ws.Range('E25').Value = '  +0.00%  '
$ws.Range('D26').Value = '4.30'
$ws.Range('E26').Value = '  -1.91%  '
$ws.Range('E27').Value = '  -2.48%  '
$ws.Range('D28').Value = '2.788.87'
$ws.Range('E28').Value = '  +1.24%  '
$ws.Range('D29').Value = '0.999'
$ws.Range('E29').Value = '  -0.02%  '
$ws.Range('D30').Value = '0.0₃0968'
$ws.Range('E30').Value = '  +1.54%  '
$ws.Range('D31').Value = '8.09'
$ws.Range('D32').Value = '498.69'
$ws.Range('E32').Value = '  -3.77%  '
$ws.Range('E33').Value = '  -2.37%  '
$ws.Range('E34').Value = '  -0.50%  '
$ws.Range('D35').Value = '0.999'
$ws.Range('E35').Value = '  -0.03%  '
$ws.Range('D36').Value = '162.90'
$ws.Range('E36').Value = '  -1.57%  '
$ws.Range('D38').Value = '18.93'
$ws.Range('E38').Value = '  -0.65%  '
$ws.Range('E39').Value = '  +0.18%  '
$ws.Range('E40').Value = '  -1.73%  '
$ws.Range('E41').Value = '  -4.25%  '
$ws.Range('E42').Value = '  -0.01%  '
$ws.Range('E43').Value = '  -1.81%  '
$ws.Range('E44').Value = '  +0.20%  '
$ws.Range('D45').Value = '0.332'
$ws.Range('E45').Value = '  +0.07%  '
$ws.Range('E46').Value = '  -0.12%  '
$ws.Range('D47').Value = '155.30'
$ws.Range('E47').Value = '  +3.81%  '
$ws.Range('D48').Value = '3.74'
$ws.Range('E48').Value = '  +0.94%  '
$ws.Range('E49').Value = '  +2.43%  '
$ws.Range('E50').Value = '  +1.90%  '
$ws.Range('D51').Value = '0.0755'
$ws.Range('E51').Value = '  -1.35%  '
